$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear removed "Variables por persona" rows (text removed from sharedStrings) ---
# A5: CONYUVIVE -> (blank), keep existing style s="8"
$ws.Range("A5").Value = ""

# A6: HIJOSDE -> (blank), keep existing style s="9"
$ws.Range("A6").Value = ""

# A9: FECHANTO -> EDAD ; style changes from the red-highlight (s=8) to the
# normal green style (s=5) used by the other rows in that band.
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "EDAD"

# A17: text stays "BUSCANDO" but style changes from s=8 to s=5
$ws.Range("A4").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# A19: text stays "INGRESOS" but style changes from s=8 to s=6
$ws.Range("A7").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# B14: THOGAR -> (blank), keep existing style s="2"
$ws.Range("B14").Value = ""

# --- Update the view state (scrolled position / active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C16").Select()
